# Add 5 new SIP transaction rows (rows 13-17) to the Mutual Fund Database sheet,
# mirroring the formatting of the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template row (row 12) holds the standard data-row formatting (number/date/currency
# styles) that we want the new rows to inherit.
$ws.Range("A12:F12").Copy()
$ws.Range("A13:F17").PasteSpecial(-4122)

# Row 13: S.No 12 - Parag Parikh Flexi cap SIP on 2025-02-03
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "Parag Parikh Flexi cap"
$ws.Range("C13").Value = "0P0000YWL1.BO"
$ws.Range("D13").Value = "2025-02-03"
$ws.Range("E13").Value = "SIP"
$ws.Range("F13").Value = 1500

# Row 14: S.No 13 - Nippon India Small cap SIP on 2025-02-03
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "Nippon India Small cap"
$ws.Range("C14").Value = "0P0000XVFY.BO"
$ws.Range("D14").Value = "2025-02-03"
$ws.Range("E14").Value = "SIP"
$ws.Range("F14").Value = 1000

# Row 15: S.No 14 - DSP Nifty 50 Equal Weightage SIP on 2025-02-03
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "DSP Nifty 50 Equal Weightage"
$ws.Range("C15").Value = "0P0001BOXZ.BO"
$ws.Range("D15").Value = "2025-02-03"
$ws.Range("E15").Value = "SIP"
$ws.Range("F15").Value = 500

# Row 16: S.No 15 - DSP Nifty Next 50 SIP on 2025-02-03
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "DSP Nifty Next 50"
$ws.Range("C16").Value = "0P0001FTFQ.BO"
$ws.Range("D16").Value = "2025-02-03"
$ws.Range("E16").Value = "SIP"
$ws.Range("F16").Value = 500

# Row 17: S.No 16 - Edelweiss Nifty 100 Quality 30 SIP on 2025-02-05
# (new fund, so its symbol/name strings are new shared strings; set the symbol (C)
# before the name (B) to match the shared-string insertion order)
$ws.Range("A17").Value = 16
$ws.Range("C17").Value = "0P0001NI59.BO"
$ws.Range("B17").Value = "Edelweiss Nifty 100 Quality 30"
$ws.Range("D17").Value = "2025-02-05"
$ws.Range("E17").Value = "SIP"
$ws.Range("F17").Value = 500

# Match the author's final cursor/selection position on the sheet.
$ws.Range("B12").Select()
$excel.CutCopyMode = $false

Write-Host "Added rows 13-17 to Mutual Fund Database sheet"
